$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.127002239227295
$ws.Range("E2").Value = 985.8829405616561
$ws.Range("F2").Value = 0.03969696128134963
$ws.Range("G2").Value = 0.03401975340484398
$ws.Range("H2").Value = 0.02989338540883184
$ws.Range("I2").Value = 0.02650547962996064
$ws.Range("J2").Value = 0.02457350479566378
$ws.Range("K2").Value = 0.02390540187109651
$ws.Range("L2").Value = 0.02236365303176127
$ws.Range("M2").Value = 0.02158710463625446
$ws.Range("N2").Value = 0.02103448310799276
$ws.Range("O2").Value = 0.02042708402507967
$ws.Range("P2").Value = 0.02009963610571553
$ws.Range("Q2").Value = 0.01981964456032652
$ws.Range("R2").Value = 0.01967325243910106
$ws.Range("S2").Value = 0.01956398952631518
$ws.Range("T2").Value = 0.01944244484830682
$ws.Range("U2").Value = 0.01936789598534594
$ws.Range("V2").Value = 0.0193551825010869
$ws.Range("W2").Value = 0.01932175538044932
$ws.Range("X2").Value = 0.019254628323309
$ws.Range("Y2").Value = 0.01921799104408686

$ws.Range("C3").Value = 1.415026187896729
$ws.Range("E3").Value = 1013.915265384405
$ws.Range("F3").Value = 0.04032535511008972
$ws.Range("G3").Value = 0.03388733915140892
$ws.Range("H3").Value = 0.03034542672478656
$ws.Range("I3").Value = 0.02769445247307296
$ws.Range("J3").Value = 0.02568364894542545
$ws.Range("K3").Value = 0.02482453885459867
$ws.Range("L3").Value = 0.02348770285554721
$ws.Range("M3").Value = 0.02274887836130471
$ws.Range("N3").Value = 0.02198987505479343
$ws.Range("O3").Value = 0.02152066233582324
$ws.Range("P3").Value = 0.02097162971737303
$ws.Range("Q3").Value = 0.02094647104481593
$ws.Range("R3").Value = 0.02065181858565355
$ws.Range("S3").Value = 0.02033824458099496
$ws.Range("T3").Value = 0.0201020281474318
$ws.Range("U3").Value = 0.01997925748846828
$ws.Range("V3").Value = 0.01994153299697672
$ws.Range("W3").Value = 0.019872346308928
$ws.Range("X3").Value = 0.01980592672057396
$ws.Range("Y3").Value = 0.01976443012445233

$ws.Range("C4").Value = 1.041002035140991
$ws.Range("E4").Value = 979.6748383368395
$ws.Range("F4").Value = 0.03881061254419911
$ws.Range("G4").Value = 0.03411561577836814
$ws.Range("H4").Value = 0.02812917005638697
$ws.Range("I4").Value = 0.02785271414410497
$ws.Range("J4").Value = 0.02504658912375245
$ws.Range("K4").Value = 0.02354919950034498
$ws.Range("L4").Value = 0.02332228827230827
$ws.Range("M4").Value = 0.02220692149746637
$ws.Range("N4").Value = 0.02165051289253221
$ws.Range("O4").Value = 0.0210261449701857
$ws.Range("P4").Value = 0.0203865675054443
$ws.Range("Q4").Value = 0.02018300753583992
$ws.Range("R4").Value = 0.01968506000186004
$ws.Range("S4").Value = 0.01968506000186004
$ws.Range("T4").Value = 0.01953695492108131
$ws.Range("U4").Value = 0.01943105049729861
$ws.Range("V4").Value = 0.01931675587423946
$ws.Range("W4").Value = 0.01914497613596267
$ws.Range("X4").Value = 0.01914497613596267
$ws.Range("Y4").Value = 0.01909697540617621

$ws.Range("C5").Value = 1.061997175216675
$ws.Range("E5").Value = 1000.61663267893
$ws.Range("F5").Value = 0.04104404569334012
$ws.Range("G5").Value = 0.03272533045882452
$ws.Range("H5").Value = 0.02933368073304323
$ws.Range("I5").Value = 0.02707185621832879
$ws.Range("J5").Value = 0.02587792970432929
$ws.Range("K5").Value = 0.02503686503593363
$ws.Range("L5").Value = 0.023200479219436
$ws.Range("M5").Value = 0.02254535743476124
$ws.Range("N5").Value = 0.02209033424821177
$ws.Range("O5").Value = 0.02154523753347399
$ws.Range("P5").Value = 0.02124117695627828
$ws.Range("Q5").Value = 0.02051786685887327
$ws.Range("R5").Value = 0.02038000764534746
$ws.Range("S5").Value = 0.02021246601500641
$ws.Range("T5").Value = 0.01998708578856482
$ws.Range("U5").Value = 0.01981781659060941
$ws.Range("V5").Value = 0.01972391046286302
$ws.Range("W5").Value = 0.01965688503021192
$ws.Range("X5").Value = 0.01958273323166826
$ws.Range("Y5").Value = 0.01950519751810779

$ws.Range("C6").Value = 1.292985677719116
$ws.Range("E6").Value = 1001.293458640235
$ws.Range("F6").Value = 0.04032861477802605
$ws.Range("G6").Value = 0.03376323281822059
$ws.Range("H6").Value = 0.02962979461653837
$ws.Range("I6").Value = 0.02675793332938933
$ws.Range("J6").Value = 0.02448528543296291
$ws.Range("K6").Value = 0.02387037332384809
$ws.Range("L6").Value = 0.02264823811911284
$ws.Range("M6").Value = 0.02146299073695688
$ws.Range("N6").Value = 0.02075669050186326
$ws.Range("O6").Value = 0.02065637049722399
$ws.Range("P6").Value = 0.02054432274955757
$ws.Range("Q6").Value = 0.0203925082325153
$ws.Range("R6").Value = 0.0201497411788545
$ws.Range("S6").Value = 0.02001426557541406
$ws.Range("T6").Value = 0.01989965703409005
$ws.Range("U6").Value = 0.01972918201977465
$ws.Range("V6").Value = 0.01965145442336613
$ws.Range("W6").Value = 0.01955368373692266
$ws.Range("X6").Value = 0.01955368373692266
$ws.Range("Y6").Value = 0.01951839100663226

$ws.Range("C7").Value = 1.181975841522217
$ws.Range("E7").Value = 989.4438487007701
$ws.Range("F7").Value = 0.03968322551281712
$ws.Range("G7").Value = 0.03346825714322793
$ws.Range("H7").Value = 0.03008371194880819
$ws.Range("I7").Value = 0.02682484295747948
$ws.Range("J7").Value = 0.02535746759807733
$ws.Range("K7").Value = 0.02345338097423296
$ws.Range("L7").Value = 0.02256340043935648
$ws.Range("M7").Value = 0.02168904992576079
$ws.Range("N7").Value = 0.02097170309220812
$ws.Range("O7").Value = 0.02083563057740072
$ws.Range("P7").Value = 0.02042853611714845
$ws.Range("Q7").Value = 0.02007091744895419
$ws.Range("R7").Value = 0.01999997431897672
$ws.Range("S7").Value = 0.01988740403494265
$ws.Range("T7").Value = 0.01961823256148919
$ws.Range("U7").Value = 0.01955441791859682
$ws.Range("V7").Value = 0.01947593083432394
$ws.Range("W7").Value = 0.01938124661760669
$ws.Range("X7").Value = 0.01931736097766218
$ws.Range("Y7").Value = 0.01928740445810468

$ws.Range("C8").Value = 1.2149977684021
$ws.Range("E8").Value = 986.370097212488
$ws.Range("F8").Value = 0.04174793149721373
$ws.Range("G8").Value = 0.03285245352127707
$ws.Range("H8").Value = 0.02935052880679257
$ws.Range("I8").Value = 0.0263539163865702
$ws.Range("J8").Value = 0.02421945397804682
$ws.Range("K8").Value = 0.02278446127313631
$ws.Range("L8").Value = 0.02244108419324938
$ws.Range("M8").Value = 0.02219732319866414
$ws.Range("N8").Value = 0.0216080221468091
$ws.Range("O8").Value = 0.02106008850045317
$ws.Range("P8").Value = 0.02045124275838433
$ws.Range("Q8").Value = 0.0203322349552415
$ws.Range("R8").Value = 0.02004573571626416
$ws.Range("S8").Value = 0.01978224546020592
$ws.Range("T8").Value = 0.01958633834880578
$ws.Range("U8").Value = 0.01950516378817512
$ws.Range("V8").Value = 0.01936582754866368
$ws.Range("W8").Value = 0.01926845312489494
$ws.Range("X8").Value = 0.01925383469158044
$ws.Range("Y8").Value = 0.01922748727509723

$ws.Range("C9").Value = 1.218042850494385
$ws.Range("E9").Value = 989.7979352131606
$ws.Range("F9").Value = 0.04001055995749826
$ws.Range("G9").Value = 0.03364266510418969
$ws.Range("H9").Value = 0.02901757713427336
$ws.Range("I9").Value = 0.0269661440034435
$ws.Range("J9").Value = 0.02452092044445193
$ws.Range("K9").Value = 0.02333670132287987
$ws.Range("L9").Value = 0.02207766911025558
$ws.Range("M9").Value = 0.02155335795214891
$ws.Range("N9").Value = 0.02104000831757854
$ws.Range("O9").Value = 0.02079973669664029
$ws.Range("P9").Value = 0.02047784163912649
$ws.Range("Q9").Value = 0.02032332700462483
$ws.Range("R9").Value = 0.0201842244155611
$ws.Range("S9").Value = 0.01997375946855839
$ws.Range("T9").Value = 0.0196723753602993
$ws.Range("U9").Value = 0.01958024164839501
$ws.Range("V9").Value = 0.01944107773493311
$ws.Range("W9").Value = 0.01936036923920216
$ws.Range("X9").Value = 0.01931113128457245
$ws.Range("Y9").Value = 0.01929430672930137

$ws.Range("C10").Value = 1.20896053314209
$ws.Range("E10").Value = 984.4815462224342
$ws.Range("F10").Value = 0.04025615495786938
$ws.Range("G10").Value = 0.03429301675504807
$ws.Range("H10").Value = 0.03082658376381462
$ws.Range("I10").Value = 0.02828490031706784
$ws.Range("J10").Value = 0.02570282249671375
$ws.Range("K10").Value = 0.02450341871890199
$ws.Range("L10").Value = 0.02306546923632629
$ws.Range("M10").Value = 0.02244590991548217
$ws.Range("N10").Value = 0.02151823358891152
$ws.Range("O10").Value = 0.02085699246341139
$ws.Range("P10").Value = 0.02055315822850005
$ws.Range("Q10").Value = 0.02026570501564194
$ws.Range("R10").Value = 0.01979698254053091
$ws.Range("S10").Value = 0.01969860948853589
$ws.Range("T10").Value = 0.01956148205167468
$ws.Range("U10").Value = 0.01944850293428409
$ws.Range("V10").Value = 0.01934511578798931
$ws.Range("W10").Value = 0.01927560198324535
$ws.Range("X10").Value = 0.01924996479213307
$ws.Range("Y10").Value = 0.01919067341564199

$ws.Range("C11").Value = 1.260037899017334
$ws.Range("E11").Value = 1004.306992222186
$ws.Range("F11").Value = 0.04169666181192185
$ws.Range("G11").Value = 0.03338760792510541
$ws.Range("H11").Value = 0.02979479767818656
$ws.Range("I11").Value = 0.02684148562643207
$ws.Range("J11").Value = 0.02516895329500047
$ws.Range("K11").Value = 0.02438405150638306
$ws.Range("L11").Value = 0.0236426744545179
$ws.Range("M11").Value = 0.02198416112209743
$ws.Range("N11").Value = 0.0216015318290662
$ws.Range("O11").Value = 0.02127513101388412
$ws.Range("P11").Value = 0.02093211773706319
$ws.Range("Q11").Value = 0.02057922870695279
$ws.Range("R11").Value = 0.02036657310340947
$ws.Range("S11").Value = 0.02009966953595505
$ws.Range("T11").Value = 0.02001610919586396
$ws.Range("U11").Value = 0.01981686077639461
$ws.Range("V11").Value = 0.01970574792200365
$ws.Range("W11").Value = 0.01969576470663202
$ws.Range("X11").Value = 0.0196065997643732
$ws.Range("Y11").Value = 0.01957713435130966
